$d = $word.ActiveDocument

$replacements = @(
  @("2024-02-08 Thursday", "2024-02-09 Friday"),
  @("30×16=", "25×64="),
  @("16×66=", "50×21="),
  @("96×77=", "80×47="),
  @("77×27=", "39×80="),
  @("50×73=", "39×91="),
  @("81×11=", "94×81="),
  @("72×82=", "63×50="),
  @("87×66=", "39×71="),
  @("86×25=", "86×24="),
  @("28×23=", "39×89="),
  @("11×66=", "47×86="),
  @("66×98=", "90×24="),
  @("54×65=", "51×54="),
  @("75×47=", "20×77="),
  @("66×94=", "99×42="),
  @("45×61=", "91×12="),
  @("32×79=", "57×98="),
  @("15×61=", "47×93="),
  @("25×29=", "48×17="),
  @("35×26=", "47×90="),
  @("19×16=", "91×12="),
  @("52×58=", "30×86="),
  @("59×38=", "85×46="),
  @("36×27=", "94×99="),
  @("39×46=", "60×97=")
)

foreach ($pair in $replacements) {
  $old = $pair[0]
  $new = $pair[1]
  $range = $d.Content
  $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
